$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (max/denominator row): add two new assignment columns ---
$ws.Range("L2").Value = 22
$ws.Range("M2").Value = 0

# --- Row 12: student now completed assignments 6 and 8 ---
$ws.Range("I12").Value = 20
$ws.Range("K12").Value = 9

# --- Row 18: student now completed assignment 6 ---
$ws.Range("I18").Value = 20

# --- Row 39: student now completed assignment 6 ---
$ws.Range("I39").Value = 20

# --- Row 43: student now has grades recorded for assignments 1-3 ---
$ws.Range("D43").Value = 41
$ws.Range("E43").Value = 18
$ws.Range("F43").Value = 24

# --- Row 11: student grades recorded for assignments 3,4,5,6,7,8; row highlighted ---
$ws.Range("F11").Value = 22
$ws.Range("G11").Value = 13
$ws.Range("H11").Value = 10
$ws.Range("I11").Value = 20
$ws.Range("J11").Value = 20
$ws.Range("K11").Value = 4

# Apply "No Fill" formatting explicitly to each populated cell in the row
# (mirrors Excel's behaviour of only stamping a style on cells that exist)
$ws.Range("A11").Interior.ColorIndex = -4142
$ws.Range("B11").Interior.ColorIndex = -4142
$ws.Range("D11").Interior.ColorIndex = -4142
$ws.Range("F11").Interior.ColorIndex = -4142
$ws.Range("G11").Interior.ColorIndex = -4142
$ws.Range("H11").Interior.ColorIndex = -4142
$ws.Range("I11").Interior.ColorIndex = -4142
$ws.Range("J11").Interior.ColorIndex = -4142
$ws.Range("K11").Interior.ColorIndex = -4142
$ws.Range("T11").Interior.ColorIndex = -4142
$ws.Range("U11").Interior.ColorIndex = -4142

# --- Update selection to match the author's last active cell ---
$ws.Range("G8").Select()
